$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.122.87"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.655.71"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "'218.28"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "'0.5299"
$ws.Range("E6").Value = "  +1.55%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "'0.2611"
$ws.Range("E8").Value = "  -2.23%  "

$ws.Range("D9").Value = "'0.06332"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "'20.38"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").Value = "'0.07751"

$ws.Range("D12").Value = "'4.502"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("D13").Value = "1.655.20"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "0.0₅8140"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "'65.23"

$ws.Range("D17").Value = "26.131.63"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "'4.537"
$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").Value = "'193.73"
$ws.Range("E20").Value = "  -0.59%  "

$ws.Range("D21").Value = "'10.03"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").Value = "'6.000"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").Value = "'140.53"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").Value = "'7.275"
$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").Value = "'16.18"
$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").Value = "'0.05943"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "'1.277"
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").Value = "'3.513"
$ws.Range("E31").Value = "  -3.52%  "

$ws.Range("D32").Value = "'3.235"
$ws.Range("E32").Value = "  -2.10%  "

$ws.Range("D33").Value = "'1.545"
$ws.Range("E33").Value = "  -5.20%  "

$ws.Range("D34").Value = "'2.412"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").Value = "'0.9458"
$ws.Range("E35").Value = "  -3.31%  "

$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("D37").Value = "'0.5633"
$ws.Range("E37").Value = "  -4.29%  "

$ws.Range("D38").Value = "'0.01608"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("D39").Value = "'5.846"
$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("D40").Value = "'0.8471"
$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "'101.28"
$ws.Range("E42").Value = "  +1.59%  "

$ws.Range("D43").Value = "1.012.02"
$ws.Range("E43").Value = "  -1.52%  "

$ws.Range("D44").Value = "1.801.17"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "'56.93"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -6.26%  "

$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").Value = "'0.4289"
$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("D49").Value = "'0.05152"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").Value = "'1.467"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").Value = "'7.738"
$ws.Range("E51").Value = "  -4.10%  "
